# UC07-Edit Request.docx — apply the commit's text edits to the Word
# object model. Each call below scopes Find/Replace to the exact
# paragraph that holds the text being changed (several snippets like
# "Crisis Management" / "main page" repeat elsewhere in the doc, so a
# document-wide Find would be ambiguous).

$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$OldText,
        [string]$NewText
    )
    $rng = $d.Paragraphs($Index).Range
    $ok = $rng.Find.Execute($OldText, $false, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed at paragraph $Index for [$OldText]"
    }
}

# Description: "Manager will change..." -> "Manager changes..."
Replace-InParagraph 5 " will change the content of a request of an incident." " changes the content of a request of an incident."

# Pre-conditions, bullet 1
Replace-InParagraph 17 "There is an active incident and is active request for that incident" "There is a request of an active request"

# Pre-conditions, bullet 2
Replace-InParagraph 18 "Crisis Management main page is open" "Crisis Board is open"

# Post-conditions, bullet 1
Replace-InParagraph 21 "System updates the request for volunteers. " "System updates the request. "

# R07-1 Main Path, step 1
Replace-InParagraph 25 "Manager searches volunteers refer to Search Volunteers use-case and selects the incident from the list of incidents" "Manager selects the incident from the crisis map and selects Resource Gathering options"

# R07-1 Main Path, step 2
Replace-InParagraph 26 "System shows information and options for the incident containing a list of requests for that incident." "System navigates to resource gathering page."

# R07-1 Main Path, step 3
$quoteOpen = [char]8220
$quoteClose = [char]8221
$oldP27 = "request and presses  " + $quoteOpen + "Edit Request" + $quoteClose + " menu"
$newP27 = "the edit request menu of a request"
Replace-InParagraph 27 $oldP27 $newP27

# R07-1 Main Path, step 5 (first line, after the <w:br/>)
$lineBreak = [char]11
$oldP29 = $lineBreak + "Edits a message for volunteers"
$newP29 = $lineBreak + "Edits the request name"
Replace-InParagraph 29 $oldP29 $newP29

# R07-1 Main Path, step 5 (second line -> becomes "Edits a message for volunteers")
Replace-InParagraph 30 "Change the selection of need list items for this request" "Edits a message for volunteers"

# R07-1 Main Path, step 5 (third line: comment anchor text + Presses 'Send Request' -> new bullet text + Presses 'Save')
$curlyOpen = [char]8216
$curlyClose = [char]8217
$oldP31 = "      " + $lineBreak + "Presses " + $curlyOpen + "Send Request" + $curlyClose + " button. "
$newP31 = "Change the selection of need list items for this request      " + $lineBreak + "Presses " + $curlyOpen + "Save" + $curlyClose + " button. "
Replace-InParagraph 31 $oldP31 $newP31

# R07-1 Main Path, step 6 (System saves...)
Replace-InParagraph 32 "System saves and starts automatic request sending and navigates to Crisis Management main page" "System saves the request and navigates to resource gathering page"

# R07-2 Alternate path
Replace-InParagraph 36 "Crisis Management main page." "resource gathering page."
